$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 updates (subject numbers)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates (CON) - B2 cleared, C2:E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0.65556633424932542
$ws.Range("D2").Value = 0.18034863288214198
$ws.Range("E2").Value = 0.96313990293282503

# Row 3 updates (STR)
$ws.Range("B3").Value = 0.12684248988640923
$ws.Range("C3").Value = 1.544713802281203
$ws.Range("D3").Value = 0.15261716424378549
$ws.Range("E3").Value = 0.92422304724005855

# Update selection to reflect the new active range
$ws.Range("B1:E3").Select()
